# Apply scheduled profit-sheet value refresh (Bahamut_Profits)
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 217171.66
$ws.Range("J17").Value = 217171.66
$ws.Range("L17").Value = 651514.98
$ws.Range("N17").Value = -651850.98
$ws.Range("H19").Value = 1740.7646
$ws.Range("I19").Value = 1383.5
$ws.Range("J19").Value = 1935.6364
$ws.Range("K19").Value = 1383.5
$ws.Range("L19").Value = 1935.6364
$ws.Range("M19").Value = -1208.5
$ws.Range("N19").Value = -2285.6364
$ws.Range("H33").Value = 160.15384
$ws.Range("I33").Value = 151
$ws.Range("K33").Value = 151
$ws.Range("M33").Value = 78
$ws.Range("H74").Value = 3398238.5
$ws.Range("I74").Value = 3640612.8
$ws.Range("K74").Value = 3640612.8
$ws.Range("M74").Value = -3639676.8
$ws.Range("H77").Value = 3398238.5
$ws.Range("I77").Value = 3640612.8
$ws.Range("K77").Value = 18203064
$ws.Range("M77").Value = -18198384
$ws.Range("H129").Value = 1014.5
$ws.Range("J129").Value = 1331.8572
$ws.Range("L129").Value = 3995.5716
$ws.Range("N129").Value = -13995.5716
$ws.Range("H137").Value = 978.2
$ws.Range("I137").Value = 714.2778
$ws.Range("K137").Value = 2142.8334
$ws.Range("M137").Value = 407.1666
$ws.Range("H139").Value = 53568.332
$ws.Range("J139").Value = 53568.332
$ws.Range("L139").Value = 53568.332
$ws.Range("N139").Value = -63848.332

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 810
$ws.Range("I61").Value = 720.4761999999999
$ws.Range("J61").Value = 1750
$ws.Range("K61").Value = 720.4761999999999
$ws.Range("L61").Value = 1750
$ws.Range("M61").Value = -508.4761999999999
$ws.Range("N61").Value = -2174
$ws.Range("H74").Value = 1144.3636
$ws.Range("I74").Value = 1109.7778
$ws.Range("K74").Value = 1109.7778
$ws.Range("M74").Value = -235.7778000000001
$ws.Range("H77").Value = 1144.3636
$ws.Range("I77").Value = 1109.7778
$ws.Range("K77").Value = 5548.889
$ws.Range("M77").Value = -1180.889
$ws.Range("H132").Value = 1357.119
$ws.Range("I132").Value = 881.9643
$ws.Range("J132").Value = 2307.4285
$ws.Range("K132").Value = 2645.8929
$ws.Range("L132").Value = 6922.2855
$ws.Range("M132").Value = -115.8928999999998
$ws.Range("N132").Value = -11982.2855
$ws.Range("H136").Value = 810
$ws.Range("I136").Value = 720.4761999999999
$ws.Range("J136").Value = 1750
$ws.Range("K136").Value = 2161.4286
$ws.Range("L136").Value = 5250
$ws.Range("M136").Value = 388.5714000000003
$ws.Range("N136").Value = -10350

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 67518.06
$ws.Range("I134").Value = 3115
$ws.Range("J134").Value = 288328.56
$ws.Range("K134").Value = 9345
$ws.Range("L134").Value = 864985.6799999999
$ws.Range("M134").Value = -6810
$ws.Range("N134").Value = -870055.6799999999

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2329.6155
$ws.Range("I31").Value = 2328.9167
$ws.Range("K31").Value = 2328.9167
$ws.Range("M31").Value = -2033.9167
$ws.Range("H34").Value = 2329.6155
$ws.Range("I34").Value = 2328.9167
$ws.Range("K34").Value = 2328.9167
$ws.Range("M34").Value = -2126.9167
$ws.Range("H58").Value = 5842.5
$ws.Range("I58").Value = 736.1111
$ws.Range("J58").Value = 51800
$ws.Range("K58").Value = 736.1111
$ws.Range("L58").Value = 51800
$ws.Range("M58").Value = -533.1111
$ws.Range("N58").Value = -52206
$ws.Range("H132").Value = 1740.8182
$ws.Range("I132").Value = 973.3889
$ws.Range("K132").Value = 2920.1667
$ws.Range("M132").Value = -390.1667000000002
$ws.Range("H134").Value = 1935.2245
$ws.Range("I134").Value = 1503.2307
$ws.Range("K134").Value = 4509.6921
$ws.Range("M134").Value = -1974.6921
$ws.Range("H136").Value = 5842.5
$ws.Range("I136").Value = 736.1111
$ws.Range("J136").Value = 51800
$ws.Range("K136").Value = 2208.3333
$ws.Range("L136").Value = 155400
$ws.Range("M136").Value = 341.6667000000002
$ws.Range("N136").Value = -160500
$ws.Range("H140").Value = 49398.43
$ws.Range("J140").Value = 49398.43
$ws.Range("L140").Value = 49398.43
$ws.Range("N140").Value = -59758.43

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 13618.182
$ws.Range("J87").Value = 21683.334
$ws.Range("L87").Value = 65050.00199999999
$ws.Range("N87").Value = -67546.00199999999
$ws.Range("H90").Value = 13618.182
$ws.Range("J90").Value = 21683.334
$ws.Range("L90").Value = 195150.006
$ws.Range("N90").Value = -207630.006
$ws.Range("H114").Value = 1316.64
$ws.Range("I114").Value = 849.8570999999999
$ws.Range("J114").Value = 1498.1666
$ws.Range("K114").Value = 2549.5713
$ws.Range("L114").Value = 4494.4998
$ws.Range("M114").Value = 704.4287000000004
$ws.Range("N114").Value = -11002.4998
$ws.Range("H132").Value = 580
$ws.Range("J132").Value = 1000
$ws.Range("L132").Value = 9000
$ws.Range("N132").Value = -14060

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 9000
$ws.Range("J52").Value = 9000
$ws.Range("L52").Value = 9000
$ws.Range("N52").Value = -9518
$ws.Range("H70").Value = 5001.6
$ws.Range("I70").Value = 4002.2856
$ws.Range("K70").Value = 4002.2856
$ws.Range("M70").Value = -3732.2856
$ws.Range("H73").Value = 5001.6
$ws.Range("I73").Value = 4002.2856
$ws.Range("K73").Value = 4002.2856
$ws.Range("M73").Value = -3066.2856
$ws.Range("H122").Value = 2041.579
$ws.Range("I122").Value = 1927.909
$ws.Range("J122").Value = 2197.875
$ws.Range("K122").Value = 5783.727000000001
$ws.Range("L122").Value = 6593.625
$ws.Range("M122").Value = -3333.727000000001
$ws.Range("N122").Value = -11493.625
$ws.Range("H126").Value = 6013
$ws.Range("J126").Value = 3014
$ws.Range("L126").Value = 9042
$ws.Range("N126").Value = -13982
$ws.Range("H132").Value = 2789.9429
$ws.Range("I132").Value = 2562.4348
$ws.Range("K132").Value = 7687.3044
$ws.Range("M132").Value = -5157.3044

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 9337.5
$ws.Range("J34").Value = 11116.667
$ws.Range("L34").Value = 11116.667
$ws.Range("N34").Value = -11460.667
$ws.Range("H132").Value = 3906.25
$ws.Range("I132").Value = 2071.4285
$ws.Range("J132").Value = 5333.3335
$ws.Range("K132").Value = 6214.2855
$ws.Range("L132").Value = 16000.0005
$ws.Range("M132").Value = -3684.2855
$ws.Range("N132").Value = -21060.0005
$ws.Range("H136").Value = 5380.8945
$ws.Range("I136").Value = 1543.3529
$ws.Range("K136").Value = 4630.0587
$ws.Range("M136").Value = -2080.0587

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1968.5883
$ws.Range("I81").Value = 2005.5
$ws.Range("J81").Value = 1880
$ws.Range("K81").Value = 4011
$ws.Range("L81").Value = 3760
$ws.Range("M81").Value = -2950
$ws.Range("N81").Value = -5882
$ws.Range("H84").Value = 1968.5883
$ws.Range("I84").Value = 2005.5
$ws.Range("J84").Value = 1880
$ws.Range("K84").Value = 20055
$ws.Range("L84").Value = 18800
$ws.Range("M84").Value = -14751
$ws.Range("N84").Value = -29408
$ws.Range("H132").Value = 1855.7778
$ws.Range("I132").Value = 1855.7778
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5567.3334
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3037.3334
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 1730.4
$ws.Range("I136").Value = 1720.48
$ws.Range("K136").Value = 5161.440000000001
$ws.Range("M136").Value = -2611.440000000001
